# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Day
$ws.Range("A2").Value = 45967

# Hourly prices B2:Z2
$ws.Range("B2").Value = 20.84
$ws.Range("C2").Value = 15.09
$ws.Range("D2").Value = 15.37
$ws.Range("E2").Value = 11.6
$ws.Range("F2").Value = 7.86
$ws.Range("G2").Value = 16.32
$ws.Range("H2").Value = 39.18
$ws.Range("I2").Value = 69.91
$ws.Range("J2").Value = 61.17
$ws.Range("K2").Value = 14.9
$ws.Range("L2").Value = 3.72
$ws.Range("M2").Value = 4.85
$ws.Range("N2").Value = 5.78
$ws.Range("O2").Value = 4.31
$ws.Range("P2").Value = 7.66
$ws.Range("Q2").Value = 7.88
$ws.Range("R2").Value = 18.33
$ws.Range("S2").Value = 75.48
$ws.Range("T2").Value = 101.16
$ws.Range("U2").Value = 119.58
$ws.Range("V2").Value = 109.62
$ws.Range("W2").Value = 92.84
$ws.Range("X2").Value = 87.56
$ws.Range("Y2").Value = 78.67
$ws.Range("Z2").Value = 41.24

# AA2 (Slot_4h_max) stays "20h-24h" - unchanged

# Slot_4h_price
$ws.Range("AB2").Value = 92.17

# Slot_2h_frist / Slot_2h_frist_price
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 110.37

# Slot_2h_second / Slot_2h_second_price
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 101.23

# Slot_min_price
$ws.Range("AG2").Value = "0h-16h"
